# Updates cryptocurrency price/volume data cells per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.738.34"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "2.826.52"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'350.65"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").Value = "'112.79"
$ws.Range("E6").Value = "  +4.53%  "

$ws.Range("E7").Value = "  +1.90%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +5.83%  "

$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").Value = "'19.98"
$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("E14").Value = "  +3.32%  "

$ws.Range("D15").Value = "3.272.04"
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("D16").Value = "'0.973"
$ws.Range("E16").Value = "  +5.43%  "

$ws.Range("D17").Value = "2.823.02"
$ws.Range("E17").Value = "  +2.02%  "

$ws.Range("D18").Value = "51.769.26"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").Value = "'3.45"
$ws.Range("E19").Value = "  +11.56%  "

$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").Value = "'13.32"
$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "'70.51"
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").Value = "'268.89"
$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("D26").Value = "'26.26"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("D29").Value = "'39.08"
$ws.Range("E29").Value = "  +7.30%  "

$ws.Range("D30").Value = "'10.57"
$ws.Range("E30").Value = "  +3.49%  "

$ws.Range("E31").Value = "  +2.04%  "

$ws.Range("D32").Value = "'52.76"
$ws.Range("E32").Value = "  +1.85%  "

$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").Value = "'0.0456"
$ws.Range("E34").Value = "  +3.20%  "

$ws.Range("E35").Value = "  +9.01%  "

$ws.Range("D36").Value = "'5.63"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "'19.05"
$ws.Range("E38").Value = "  +4.65%  "

$ws.Range("D39").Value = "'3.23"
$ws.Range("E39").Value = "  +2.86%  "

$ws.Range("E40").Value = "  +2.62%  "

$ws.Range("E41").Value = "  +1.59%  "

$ws.Range("D42").Value = "'2.53"
$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("D43").Value = "'121.65"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("D45").Value = "'22.02"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.185.34"
$ws.Range("E46").Value = "  +4.22%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.48"
$ws.Range("E47").Value = "  +7.36%  "

$ws.Range("E48").Value = "  +8.08%  "

$ws.Range("D49").Value = "'0.246"
$ws.Range("E49").Value = "  +27.85%  "

$ws.Range("D50").Value = "'0.982"
$ws.Range("E50").Value = "  +8.77%  "

$ws.Range("D51").Value = "'5.51"
$ws.Range("E51").Value = "  +1.72%  "
